# Reorder two "PROFESSIONAL EXPERIENCE" job blocks:
#  1. Move the "Senior Analyst - Myers Research" block (Heading 3 title +
#     subtitle + 3 bullet paragraphs = 5 paragraphs total) so it appears
#     immediately before the "Research Director - PCCC" block (previously
#     it came right after that block).
#  2. Move the "Field Director - The Feldman Group" block (same shape) so
#     it appears immediately before the "Programmer - Lake Research
#     Partners" block (previously it came right after that block).
#
# Word's Range.Cut()/Range.Paste() moves the paragraph text but the
# destination paragraph's own formatting can win out for the pasted-in
# paragraph mark, so each moved heading's "Heading 3" style is
# re-applied explicitly after the paste to be safe.

$d = $word.ActiveDocument

# --- Step 1: Senior Analyst - Myers Research -> before Research Director - PCCC ---

$fSenior = $d.Content.Duplicate
$fSenior.Find.Execute("Senior Analyst - Myers Research", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$seniorIdx = $fSenior.Paragraphs.Item(1).Index

$fResearch = $d.Content.Duplicate
$fResearch.Find.Execute("Research Director - PCCC", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$researchIdx = $fResearch.Paragraphs.Item(1).Index

# Block = heading paragraph + the following 4 body paragraphs.
$startP = $d.Paragraphs.Item($seniorIdx)
$endP = $d.Paragraphs.Item($seniorIdx + 4)
$moveRange = $d.Range($startP.Range.Start, $endP.Range.End)
$moveRange.Cut()

$targetP = $d.Paragraphs.Item($researchIdx)
$insertAt = $d.Range($targetP.Range.Start, $targetP.Range.Start)
$insertAt.Paste()

# Re-assert the heading style on the paragraph that was just pasted in.
$d.Paragraphs.Item($researchIdx).Style = "Heading 3"

# --- Step 2: Field Director - The Feldman Group -> before Programmer - Lake Research Partners ---

$fField = $d.Content.Duplicate
$fField.Find.Execute("Field Director - The Feldman Group", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fieldIdx = $fField.Paragraphs.Item(1).Index

$fProgrammer = $d.Content.Duplicate
$fProgrammer.Find.Execute("Programmer - Lake Research Partners", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$programmerIdx = $fProgrammer.Paragraphs.Item(1).Index

$startP2 = $d.Paragraphs.Item($fieldIdx)
$endP2 = $d.Paragraphs.Item($fieldIdx + 4)
$moveRange2 = $d.Range($startP2.Range.Start, $endP2.Range.End)
$moveRange2.Cut()

$targetP2 = $d.Paragraphs.Item($programmerIdx)
$insertAt2 = $d.Range($targetP2.Range.Start, $targetP2.Range.Start)
$insertAt2.Paste()

# Re-assert the heading style on the paragraph that was just pasted in.
$d.Paragraphs.Item($programmerIdx).Style = "Heading 3"
